$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.364.82"
$ws.Range("E2").Value = "  +1.32%  "

$ws.Range("D3").Value = "3.119.35"
$ws.Range("E3").Value = "  +1.96%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "396.12"
$ws.Range("E5").Value = "  +2.98%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.94"
$ws.Range("E6").Value = "  +0.35%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.540"
$ws.Range("E7").Value = "  -0.98%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.600"
$ws.Range("E9").Value = "  +2.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.91"
$ws.Range("E10").Value = "  +1.91%  "

$ws.Range("E11").Value = "  +0.91%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0861"
$ws.Range("E12").Value = "  -0.58%  "

$ws.Range("D13").Value = "3.609.52"
$ws.Range("E13").Value = "  +1.70%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.80"
$ws.Range("E14").Value = "  +0.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.85"
$ws.Range("E15").Value = "  +0.94%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.05"
$ws.Range("E16").Value = "  +7.18%  "

$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.41"
$ws.Range("E17").Value = "  +8.44%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.078.79"
$ws.Range("E18").Value = "  -0.22%  "

$ws.Range("D19").Value = "52.205.53"
$ws.Range("E19").Value = "  +0.90%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.24"
$ws.Range("E20").Value = "  +2.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.68"
$ws.Range("E21").Value = "  +1.57%  "

$ws.Range("E22").Value = "  +0.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.03"
$ws.Range("E23").Value = "  +1.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.79"
$ws.Range("E24").Value = "  -0.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.22"
$ws.Range("E25").Value = "  +1.68%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.07"
$ws.Range("E26").Value = "  -5.18%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.64"
$ws.Range("E27").Value = "  +2.28%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.36"
$ws.Range("E28").Value = "  +0.92%  "

$ws.Range("E29").Value = "  -2.38%  "

$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.108"
$ws.Range("E31").Value = "  -0.19%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.94"
$ws.Range("E32").Value = "  +5.86%  "

$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "36.81"
$ws.Range("E33").Value = "  +6.57%  "

$ws.Range("B34").Value = "VeChain"
$ws.Range("C34").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0491"
$ws.Range("E34").Value = "  +10.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.09"
$ws.Range("E35").Value = "  +0.95%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "49.89"
$ws.Range("E36").Value = "  -1.12%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.42"
$ws.Range("E38").Value = "  +0.75%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.293"
$ws.Range("E39").Value = "  +0.85%  "

$ws.Range("E40").Value = "  +8.70%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.69"
$ws.Range("E41").Value = "  +5.67%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.08"
$ws.Range("E42").Value = "  -0.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "130.62"
$ws.Range("E43").Value = "  +1.47%  "

$ws.Range("E44").Value = "  -0.39%  "

$ws.Range("E45").Value = "  +0.12%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.33"
$ws.Range("E46").Value = "  +1.33%  "

$ws.Range("E47").Value = "  -3.14%  "

$ws.Range("E48").Value = "  -1.37%  "

$ws.Range("D49").Value = "2.089.29"
$ws.Range("E49").Value = "  +1.91%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0527"
$ws.Range("E50").Value = "  +34.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.928"
$ws.Range("E51").Value = "  +11.26%  "
